# Codis_Territori_Continents.xlsx — rework the "Continents" reference table
# so it uses short alpha codes (AF, AM, AN, ...) instead of the old
# numeric identifiers (00001, 00002, ...), and rename the lookup columns
# from "Identificador"/"Nom Català" to "Codi"/"Nom".
#
# Commit message: "Modelatge Països, continents, subcontinents"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Continents")

# --- Table header (row 3) -------------------------------------------------
$ws.Range("A3").Value = "Codi"
$ws.Range("B3").Value = "Nom"

# --- Table body (rows 4-11): replace numeric identifiers with alpha codes -
# Continent names in column B stay the same (only "Altres/diversos" gets
# re-cased to "Altres/Diversos"); only column A's codes change.
$ws.Range("A4").Value = "AF"   # Àfrica
$ws.Range("A5").Value = "AM"   # Amèrica
$ws.Range("A6").Value = "AN"   # Antàrtida i Territoris Propers
$ws.Range("A7").Value = "AS"   # Àsia
$ws.Range("A8").Value = "EU"   # Europa
$ws.Range("A9").Value = "OC"   # Oceania
$ws.Range("A10").Value = "NC"  # No consta
$ws.Range("A11").Value = "AD"  # Altres/Diversos

$ws.Range("B11").Value = "Altres/Diversos"
